$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "76.525.12"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.037.40"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.46%  "

$ws.Range("E4").Value = "  -0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "201.03"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "630.48"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.47%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.552"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.36%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.204"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.034.37"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.43%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.437"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("E13").Value = "  +3.95%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.595.76"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.48%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "29.42"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.41%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "76.428.80"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000190"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.35%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.018.99"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.50%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.70%  "

$ws.Range("E20").Value = "  +1.93%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "376.17"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.92%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.36"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.22%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.206.47"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "73.07"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.07%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.88"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("E29").Value = "  +0.66%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("E31").Value = "  +8.33%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.41"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "514.67"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("E34").Value = "  +7.13%  "

$ws.Range("E35").Value = "  -0.01%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "20.95"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.64%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "164.22"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.46%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "20.01"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "

$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "193.14"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.50%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.383"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +11.04%  "

$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("E43").Value = "  +0.21%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.08"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.54"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.67%  "

$ws.Range("E46").Value = "  +4.22%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.66"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("E48").Value = "  +6.53%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.701"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.09%  "

$ws.Range("E51").Value = "  +4.72%  "
